# "Generate Report for Handoff"
#
# Regenerates the localization-status report: status text moves from
# "Handed back: in sync with en-US" to "Ready for handoff", the
# handoff/handback timestamps and priority are refreshed, and a new
# "version is stale" error message is recorded against the first e2e
# file for each locale.

$wb = $excel.ActiveWorkbook

$statusNew = "Ready for handoff"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8bd7834c855ed569f07a45d4a5cafd9fcb1bf596/e2e/2a8951f5-0894-405c-98a5-5925ac54b50c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aeb8bcf920a80e6e544d316ab12c29cc00800061/e2e/2a8951f5-0894-405c-98a5-5925ac54b50c.md."

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("G2").Value = "2016-10-19 17:56:31"
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Range("G3").Value = "2016-10-19 17:56:31"

# Columns E/F got narrower once the new (shorter) status text landed.
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusNew
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-10-19 17:56:19"
$wsZhCn.Range("P2").Value = $errorDetail
$wsZhCn.Range("C3").Value = $statusNew
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-10-19 17:56:19"

$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusNew
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-10-19 17:56:31"
$wsDeDe.Range("P2").Value = $errorDetail
$wsDeDe.Range("C3").Value = $statusNew
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-10-19 17:56:31"

$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
